# Insert a new data row right above the current row 164 ("Mango" price
# record), pushing the existing rows 164-240 down to 165-241, then
# populate the freshly inserted row 164 with the new weekly price entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(164).Insert()

$ws.Cells.Item(164, 1).Value  = 4
$ws.Cells.Item(164, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(164, 3).Value  = "Los Lagos"
$ws.Cells.Item(164, 4).Value  = 44806
$ws.Cells.Item(164, 5).Value  = 10
$ws.Cells.Item(164, 6).Value  = "Fruta"
$ws.Cells.Item(164, 7).Value  = 100108
$ws.Cells.Item(164, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(164, 9).Value  = 100108002
$ws.Cells.Item(164, 10).Value = "Mango"
$ws.Cells.Item(164, 11).Value = "Sin especificar"
$ws.Cells.Item(164, 12).Value = "Primera"
$ws.Cells.Item(164, 13).Value = 200
$ws.Cells.Item(164, 14).Value = 11000
$ws.Cells.Item(164, 15).Value = 12000
$ws.Cells.Item(164, 16).Value = 11500
$ws.Cells.Item(164, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(164, 18).Value = "Brasil"
$ws.Cells.Item(164, 19).Value = 2875
$ws.Cells.Item(164, 20).Value = 4
